$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 2 4 '63.108.59'
$ws.Cells.Item(2, 5).Value = '  +2.96%  '
Set-TextValue 3 4 '3.042.07'
$ws.Cells.Item(3, 5).Value = '  +1.88%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
Set-TextValue 5 4 '596.14'
$ws.Cells.Item(5, 5).Value = '  -0.65%  '
Set-TextValue 6 4 '154.54'
$ws.Cells.Item(6, 5).Value = '  +7.76%  '
$ws.Cells.Item(7, 5).Value = '  -0.06%  '
Set-TextValue 8 4 '3.042.22'
$ws.Cells.Item(8, 5).Value = '  +2.04%  '
$ws.Cells.Item(9, 5).Value = '  +0.19%  '
Set-TextValue 10 4 '6.88'
$ws.Cells.Item(10, 5).Value = '  +13.84%  '
$ws.Cells.Item(11, 5).Value = '  +4.00%  '
Set-TextValue 12 4 '0.465'
$ws.Cells.Item(12, 5).Value = '  +2.43%  '
Set-TextValue 13 4 '0.0000235'
$ws.Cells.Item(13, 5).Value = '  +3.43%  '
Set-TextValue 14 4 '35.99'
$ws.Cells.Item(14, 5).Value = '  +4.89%  '
$ws.Cells.Item(15, 5).Value = '  +2.00%  '
Set-TextValue 16 4 '3.550.18'
$ws.Cells.Item(16, 5).Value = '  +1.93%  '
Set-TextValue 17 4 '7.11'
$ws.Cells.Item(17, 5).Value = '  +2.50%  '
Set-TextValue 18 4 '63.122.68'
$ws.Cells.Item(18, 5).Value = '  +2.85%  '
Set-TextValue 19 4 '3.044.65'
$ws.Cells.Item(19, 5).Value = '  +2.11%  '
Set-TextValue 20 4 '455.23'
$ws.Cells.Item(20, 5).Value = '  +1.41%  '
Set-TextValue 21 4 '14.33'
$ws.Cells.Item(21, 5).Value = '  +2.74%  '
Set-TextValue 22 4 '0.700'
$ws.Cells.Item(22, 5).Value = '  +2.55%  '
$ws.Cells.Item(23, 5).Value = '  +3.37%  '
Set-TextValue 24 4 '83.13'
$ws.Cells.Item(24, 5).Value = '  +2.04%  '
Set-TextValue 25 4 '11.27'
$ws.Cells.Item(25, 5).Value = '  +5.50%  '
$ws.Cells.Item(26, 5).Value = '  +4.71%  '
Set-TextValue 27 4 '12.38'
$ws.Cells.Item(27, 5).Value = '  +3.95%  '
$ws.Cells.Item(28, 5).Value = '  +0.07%  '
Set-TextValue 29 4 '7.46'
$ws.Cells.Item(29, 5).Value = '  +4.26%  '
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 30 4 '2.70'
$ws.Cells.Item(30, 5).Value = '  +0.85%  '
$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 31 4 '2.24'
$ws.Cells.Item(31, 5).Value = '  +9.31%  '
$ws.Cells.Item(32, 5).Value = '  -0.15%  '
Set-TextValue 33 4 '27.75'
$ws.Cells.Item(33, 5).Value = '  +2.12%  '
$ws.Cells.Item(34, 5).Value = '  +1.06%  '
Set-TextValue 35 4 '0.0₃0862'
$ws.Cells.Item(35, 5).Value = '  +5.08%  '
$ws.Cells.Item(36, 5).Value = '  +2.89%  '
$ws.Cells.Item(37, 5).Value = '  +3.48%  '
Set-TextValue 38 4 '3.18'
$ws.Cells.Item(38, 5).Value = '  +11.55%  '
$ws.Cells.Item(39, 5).Value = '  +7.45%  '
Set-TextValue 40 4 '2.12'
$ws.Cells.Item(40, 5).Value = '  +3.13%  '
Set-TextValue 41 4 '50.46'
$ws.Cells.Item(41, 5).Value = '  +0.05%  '
Set-TextValue 42 4 '9.15'
$ws.Cells.Item(42, 5).Value = '  -0.50%  '
$ws.Cells.Item(43, 5).Value = '  +13.09%  '
Set-TextValue 44 4 '43.79'
$ws.Cells.Item(44, 5).Value = '  +10.48%  '
Set-TextValue 45 4 '394.13'
$ws.Cells.Item(45, 5).Value = '  -0.52%  '
$ws.Cells.Item(46, 5).Value = '  +3.07%  '
Set-TextValue 47 4 '2.726.14'
$ws.Cells.Item(47, 5).Value = '  +1.56%  '
Set-TextValue 48 4 '132.26'
$ws.Cells.Item(48, 5).Value = '  +1.12%  '
$ws.Cells.Item(49, 5).Value = '  +0.02%  '
$ws.Cells.Item(50, 5).Value = '  +6.95%  '
Set-TextValue 51 4 '24.58'
$ws.Cells.Item(51, 5).Value = '  +5.05%  '
